$d = $word.ActiveDocument

# 1. Update the letter date.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# 2. Split "944 Deer Meadow Ct, San Jose CA 95122" into two paragraphs:
#    "944 Deer Meadow Ct" and a new paragraph "San Jose, CA 95122".
$addrPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "944 Deer Meadow Ct, San Jose CA 95122*") {
        $addrPara = $p
        break
    }
}
$addrRange = $addrPara.Range
$addrRange.Text = "944 Deer Meadow Ct"
$newRange = $addrRange.InsertParagraphAfter()
$afterPara = $addrPara.Next()
$afterPara.Range.Text = "San Jose, CA 95122"

# 3. Remove the empty "No Spacing" paragraph that follows
#    "San Jose - Deerfield Homeowners Association Board of Directors".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Style.NameLocal -eq "No Spacing") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text -like "*Board of Directors*") {
            $p.Range.Delete() | Out-Null
            break
        }
    }
}
